$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1076
$ws.Range("F8").Value = 583
$ws.Range("F9").Value = 1507
$ws.Range("F11").Value = 1405
$ws.Range("F12").Value = 3046
$ws.Range("F13").Value = 550
$ws.Range("F14").Value = 1709
$ws.Range("F15").Value = 1777
$ws.Range("F17").Value = 258
$ws.Range("F21").Value = 1168
$ws.Range("F22").Value = 382
$ws.Range("F23").Value = 426
$ws.Range("F25").Value = 4181
$ws.Range("F26").Value = 725
$ws.Range("F28").Value = 1605
$ws.Range("F29").Value = 9
$ws.Range("F30").Value = 72

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 18
$ws.Range("F13").Value = 94

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 803

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 803
$ws.Range("F9").Value = 18
$ws.Range("F16").Value = 1076
$ws.Range("F19").Value = 583
$ws.Range("F20").Value = 1507
$ws.Range("F22").Value = 1405
$ws.Range("F23").Value = 3046
$ws.Range("F24").Value = 550
$ws.Range("F25").Value = 1709
$ws.Range("F26").Value = 1777
$ws.Range("F28").Value = 258
$ws.Range("F34").Value = 1168
$ws.Range("F35").Value = 382
$ws.Range("F36").Value = 426
$ws.Range("F38").Value = 4181
$ws.Range("F39").Value = 725
$ws.Range("F41").Value = 1605
$ws.Range("F42").Value = 94
$ws.Range("F44").Value = 9
$ws.Range("F45").Value = 72
